# Scheduled market-data refresh: update cached Leve profit figures
# (currentAveragePrice*/LevePrice*/LeveProfit* columns H-N) across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Cells.Item(11, 8).Value = 26.7
$ws.Cells.Item(11, 9).Value = 26.7
$ws.Cells.Item(11, 11).Value = 26.7
$ws.Cells.Item(11, 13).Value = 113.3
# Row 19
$ws.Cells.Item(19, 8).Value = 1035.4193
$ws.Cells.Item(19, 9).Value = 592.35
$ws.Cells.Item(19, 10).Value = 1841
$ws.Cells.Item(19, 11).Value = 592.35
$ws.Cells.Item(19, 12).Value = 1841
$ws.Cells.Item(19, 13).Value = -417.35
$ws.Cells.Item(19, 14).Value = -2191
# Row 88
$ws.Cells.Item(88, 8).Value = 6079.1177
$ws.Cells.Item(88, 9).Value = 4278.8335
$ws.Cells.Item(88, 10).Value = 7061.091
$ws.Cells.Item(88, 11).Value = 4278.8335
$ws.Cells.Item(88, 12).Value = 7061.091
$ws.Cells.Item(88, 13).Value = -3872.8335
$ws.Cells.Item(88, 14).Value = -7873.091
# Row 91
$ws.Cells.Item(91, 8).Value = 6079.1177
$ws.Cells.Item(91, 9).Value = 4278.8335
$ws.Cells.Item(91, 10).Value = 7061.091
$ws.Cells.Item(91, 11).Value = 4278.8335
$ws.Cells.Item(91, 12).Value = 7061.091
$ws.Cells.Item(91, 13).Value = -2874.8335
$ws.Cells.Item(91, 14).Value = -9869.091
# Row 129
$ws.Cells.Item(129, 8).Value = 773
$ws.Cells.Item(129, 10).Value = 1003.4
$ws.Cells.Item(129, 12).Value = 3010.2
$ws.Cells.Item(129, 14).Value = -13010.2
# Row 137
$ws.Cells.Item(137, 8).Value = 1856.25
$ws.Cells.Item(137, 9).Value = 1335.7142
$ws.Cells.Item(137, 10).Value = 2585
$ws.Cells.Item(137, 11).Value = 4007.1426
$ws.Cells.Item(137, 12).Value = 7755
$ws.Cells.Item(137, 13).Value = -1457.1426
$ws.Cells.Item(137, 14).Value = -12855

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 2526.1667
$ws.Cells.Item(2, 9).Value = 2833.5
$ws.Cells.Item(2, 10).Value = 2218.8333
$ws.Cells.Item(2, 11).Value = 2833.5
$ws.Cells.Item(2, 12).Value = 2218.8333
$ws.Cells.Item(2, 13).Value = -2720.5
$ws.Cells.Item(2, 14).Value = -2444.8333
# Row 32
$ws.Cells.Item(32, 8).Value = 5433.171
$ws.Cells.Item(32, 9).Value = 3602.5781
$ws.Cells.Item(32, 10).Value = 15196.333
$ws.Cells.Item(32, 11).Value = 3602.5781
$ws.Cells.Item(32, 12).Value = 15196.333
$ws.Cells.Item(32, 13).Value = -3315.5781
$ws.Cells.Item(32, 14).Value = -15770.333
# Row 63
$ws.Cells.Item(63, 8).Value = 3103.4
$ws.Cells.Item(63, 9).Value = 2650.2
$ws.Cells.Item(63, 10).Value = 3330
$ws.Cells.Item(63, 11).Value = 2650.2
$ws.Cells.Item(63, 12).Value = 3330
$ws.Cells.Item(63, 13).Value = -1964.2
$ws.Cells.Item(63, 14).Value = -4702
# Row 66
$ws.Cells.Item(66, 8).Value = 3103.4
$ws.Cells.Item(66, 9).Value = 2650.2
$ws.Cells.Item(66, 10).Value = 3330
$ws.Cells.Item(66, 11).Value = 13251
$ws.Cells.Item(66, 12).Value = 16650
$ws.Cells.Item(66, 13).Value = -9819
$ws.Cells.Item(66, 14).Value = -23514
# Row 116
$ws.Cells.Item(116, 8).Value = 2526.1667
$ws.Cells.Item(116, 9).Value = 2833.5
$ws.Cells.Item(116, 10).Value = 2218.8333
$ws.Cells.Item(116, 11).Value = 2833.5
$ws.Cells.Item(116, 12).Value = 2218.8333
$ws.Cells.Item(116, 13).Value = -539.5
$ws.Cells.Item(116, 14).Value = -6806.8333

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 2526.1667
$ws.Cells.Item(3, 9).Value = 2833.5
$ws.Cells.Item(3, 10).Value = 2218.8333
$ws.Cells.Item(3, 11).Value = 2833.5
$ws.Cells.Item(3, 12).Value = 2218.8333
$ws.Cells.Item(3, 13).Value = -2719.5
$ws.Cells.Item(3, 14).Value = -2446.8333

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 1321.0135
$ws.Cells.Item(31, 9).Value = 831.6667
$ws.Cells.Item(31, 10).Value = 2080.3447
$ws.Cells.Item(31, 11).Value = 831.6667
$ws.Cells.Item(31, 12).Value = 2080.3447
$ws.Cells.Item(31, 13).Value = -536.6667
$ws.Cells.Item(31, 14).Value = -2670.3447
# Row 34
$ws.Cells.Item(34, 8).Value = 1321.0135
$ws.Cells.Item(34, 9).Value = 831.6667
$ws.Cells.Item(34, 10).Value = 2080.3447
$ws.Cells.Item(34, 11).Value = 831.6667
$ws.Cells.Item(34, 12).Value = 2080.3447
$ws.Cells.Item(34, 13).Value = -629.6667
$ws.Cells.Item(34, 14).Value = -2484.3447
# Row 58
$ws.Cells.Item(58, 8).Value = 5919.5
$ws.Cells.Item(58, 9).Value = 6569
$ws.Cells.Item(58, 10).Value = 1806
$ws.Cells.Item(58, 11).Value = 6569
$ws.Cells.Item(58, 12).Value = 1806
$ws.Cells.Item(58, 13).Value = -6366
$ws.Cells.Item(58, 14).Value = -2212
# Row 62
$ws.Cells.Item(62, 8).Value = 2321.8262
$ws.Cells.Item(62, 9).Value = 2293.75
$ws.Cells.Item(62, 10).Value = 2386
$ws.Cells.Item(62, 11).Value = 2293.75
$ws.Cells.Item(62, 12).Value = 2386
$ws.Cells.Item(62, 13).Value = -1669.75
$ws.Cells.Item(62, 14).Value = -3634
# Row 65
$ws.Cells.Item(65, 8).Value = 2321.8262
$ws.Cells.Item(65, 9).Value = 2293.75
$ws.Cells.Item(65, 10).Value = 2386
$ws.Cells.Item(65, 11).Value = 11468.75
$ws.Cells.Item(65, 12).Value = 11930
$ws.Cells.Item(65, 13).Value = -8348.75
$ws.Cells.Item(65, 14).Value = -18170
# Row 136
$ws.Cells.Item(136, 8).Value = 5919.5
$ws.Cells.Item(136, 9).Value = 6569
$ws.Cells.Item(136, 10).Value = 1806
$ws.Cells.Item(136, 11).Value = 19707
$ws.Cells.Item(136, 12).Value = 5418
$ws.Cells.Item(136, 13).Value = -17157
$ws.Cells.Item(136, 14).Value = -10518

$ws = $wb.Worksheets.Item("CUL")
# Row 63
$ws.Cells.Item(63, 8).Value = 5606.357
$ws.Cells.Item(63, 9).Value = 3954
$ws.Cells.Item(63, 10).Value = 6057
$ws.Cells.Item(63, 11).Value = 11862
$ws.Cells.Item(63, 12).Value = 18171
$ws.Cells.Item(63, 13).Value = -11113
$ws.Cells.Item(63, 14).Value = -19669
# Row 66
$ws.Cells.Item(66, 8).Value = 5606.357
$ws.Cells.Item(66, 9).Value = 3954
$ws.Cells.Item(66, 10).Value = 6057
$ws.Cells.Item(66, 11).Value = 35586
$ws.Cells.Item(66, 12).Value = 54513
$ws.Cells.Item(66, 13).Value = -31842
$ws.Cells.Item(66, 14).Value = -62001
# Row 70
$ws.Cells.Item(70, 8).Value = 5941.1
$ws.Cells.Item(70, 9).Value = 4004
$ws.Cells.Item(70, 10).Value = 6771.2856
$ws.Cells.Item(70, 11).Value = 12012
$ws.Cells.Item(70, 12).Value = 20313.8568
$ws.Cells.Item(70, 13).Value = -11697
$ws.Cells.Item(70, 14).Value = -20943.8568
# Row 73
$ws.Cells.Item(73, 8).Value = 5941.1
$ws.Cells.Item(73, 9).Value = 4004
$ws.Cells.Item(73, 10).Value = 6771.2856
$ws.Cells.Item(73, 11).Value = 12012
$ws.Cells.Item(73, 12).Value = 20313.8568
$ws.Cells.Item(73, 13).Value = -10920
$ws.Cells.Item(73, 14).Value = -22497.8568
# Row 131
$ws.Cells.Item(131, 8).Value = 914.3
$ws.Cells.Item(131, 10).Value = 918.1818
$ws.Cells.Item(131, 12).Value = 2754.5454
$ws.Cells.Item(131, 14).Value = -12834.5454
# Row 137
$ws.Cells.Item(137, 8).Value = 19341634
$ws.Cells.Item(137, 9).Value = 5243.3335
$ws.Cells.Item(137, 10).Value = 21154422
$ws.Cells.Item(137, 11).Value = 15730.0005
$ws.Cells.Item(137, 12).Value = 63463266
$ws.Cells.Item(137, 13).Value = -10630.0005
$ws.Cells.Item(137, 14).Value = -63473466

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81, 8).Value = 1206.3334
$ws.Cells.Item(81, 10).Value = 1233.3334
$ws.Cells.Item(81, 12).Value = 2466.6668
$ws.Cells.Item(81, 14).Value = -4588.6668
# Row 84
$ws.Cells.Item(84, 8).Value = 1206.3334
$ws.Cells.Item(84, 10).Value = 1233.3334
$ws.Cells.Item(84, 12).Value = 12333.334
$ws.Cells.Item(84, 14).Value = -22941.334
# Row 132
$ws.Cells.Item(132, 8).Value = 2411.0754
$ws.Cells.Item(132, 9).Value = 2721.353
$ws.Cells.Item(132, 10).Value = 1855.8422
$ws.Cells.Item(132, 11).Value = 8164.059
$ws.Cells.Item(132, 12).Value = 5567.5266
$ws.Cells.Item(132, 13).Value = -5634.059
$ws.Cells.Item(132, 14).Value = -10627.5266
